$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new label "A" in H14 (new column used for an extra result label)
$ws.Range("H14").Value = "A"

# Add a new label/value pair in row 16: "FlujoCentro" = 0.02
$ws.Range("F16").Value = "FlujoCentro"
$ws.Range("G16").Value = 0.02
